# Updated function in Excel to handle Uppercase labels
#
# Replaces the RemoveSpecial(PROPER(...)) helper usage with a new
# fLetter() VBA-style helper (documented on the "Lookup" sheet) that
# upper-cases the first letter of every word, and points the "For field"
# example formula at RemoveSpecial(fLetter(E62)) instead of
# RemoveSpecial(PROPER(E62)).

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("Data")
$lookup = $wb.Worksheets.Item("Lookup")

# --- Data sheet: just a selection change (F4 -> B4) -----------------------
$data.Range("B4").Select()

# --- Lookup sheet: make room for the new fLetter() function body ----------
# Existing rows 13:25 (the "For field" / "For EDT" example block) shift
# down to rows 22:34, exactly like inserting 9 blank rows above row 13.
$lookup.Rows("12:20").Insert()

# RemoveSpecial's header (H3) becomes bold, matching the other function
# header/section labels on the sheet.
$lookup.Range("H3").Font.Bold = $true

# --- Write the new fLetter() function source into H12:H20 -----------------
$lookup.Range("H12").Value = "Function fLetter(str As String) As String"
$lookup.Range("H12").Font.Bold = $true

$lookup.Range("H13").Value = "    Dim strarr() As String"
$lookup.Range("H14").Value = "    strarr = Split(str)"
$lookup.Range("H15").Value = "    Dim i As Long"
$lookup.Range("H16").Value = "    For i = LBound(strarr) To UBound(strarr)"
$lookup.Range("H17").Value = "        strarr(i) = UCase`$(Left`$(strarr(i), 1)) & Mid`$(strarr(i), 2)"
$lookup.Range("H18").Value = "    Next i"
$lookup.Range("H19").Value = "    fLetter = Join(strarr, "" "")"
$lookup.Range("H20").Value = "End Function"

# --- Fix up the "For field" example formula (now at H23) ------------------
# It used to read =RemoveSpecial(PROPER(E62)); it now reads
# =RemoveSpecial(fLetter(E62)). Cells in this block are formatted as Text
# ("@") so the formula-looking string is stored as literal text, not an
# actual formula. We stage the new literal text on a scratch cell (using a
# leading apostrophe to force text) and then copy/paste-special (values) it
# onto H23 so the destination keeps its plain "Text" style instead of
# picking up a quote-prefix style variant.
$scratch = $lookup.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "'=RemoveSpecial(fLetter(E62))"
$scratch.Copy()
$lookup.Range("H23").PasteSpecial(-4163)
$scratch.Clear()

$lookup.Range("H24").Select()

# Restore "Data" as the active sheet/tab (it must stay tabSelected="1"),
# while the "Lookup" sheet keeps its own H24 selection set above.
$data.Activate()
$data.Range("B4").Select()

$excel.CutCopyMode = 0
